# Weekly price-data update: insert a new week (fecha 45075) of "Ajo" quotes
# right after the existing row 1086, pushing the remaining historical rows
# down by 7 (old 1087-1129 become 1094-1136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows starting at row 1087 (this shifts rows 1087:1129 -> 1094:1136)
$ws.Range("1087:1093").Insert()

$newRows = @(
    @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",45075,13,100112003,"Ajo","Chino","Extra",500,18000,18500,18200,'$/malla 10 kilos',"China",1820,10,"Hortaliza"),
    @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",45075,13,100112003,"Ajo","Chino","Primera",2600,14500,15000,14788,'$/caja 10 kilos',"China",1479,10,"Hortaliza"),
    @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",45075,13,100112003,"Ajo","Chino","Primera",800,15000,16000,15625,'$/malla 10 kilos',"China",1562,10,"Hortaliza"),
    @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",45075,13,100112003,"Ajo","Rosado","1a (guarda)",800,6000,6000,6000,'$/trenza 50 unidades',"Provincia de Talagante",1200,5,"Hortaliza"),
    @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",45075,13,100112003,"Ajo","Rosado","2a (guarda)",500,4500,4500,4500,'$/trenza 50 unidades',"Provincia de Talagante",900,5,"Hortaliza"),
    @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",45075,13,100112003,"Ajo","Rosado","3a (guarda)",400,2500,2500,2500,'$/trenza 50 unidades',"Provincia de Talagante",500,5,"Hortaliza"),
    @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",45075,13,100112003,"Ajo","Rosado","Extra Guarda",400,7000,7000,7000,'$/trenza 50 unidades',"Provincia de Talagante",1400,5,"Hortaliza")
)

$startRow = 1087
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Output "Inserted 7 rows of new Ajo price data at rows 1087-1093"
